$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Copy the date-formatted style from an existing row onto the new rows'
# date columns (C and D) before writing values, so the existing shared
# date style is reused instead of a new numFmt being created.
$ws.Range("C146:D146").Copy()
$ws.Range("C147:D148").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 147: Seinfeldia
$ws.Range("A147").Value = "Seinfeldia"
$ws.Range("B147").Value = "Jennifer Keishin Armstrong"
$ws.Range("C147").Value = (Get-Date -Year 2020 -Month 11 -Day 23 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D147").Value = (Get-Date -Year 2020 -Month 11 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E147").Value = "seinfeld;tv;history;comedy;writing"
$ws.Range("F147").Value = "Audio"
$ws.Range("G147").Value = "9 Hours 59 Mins"
$ws.Range("H147").Value = 3
$ws.Range("I147").Value = $false

# Row 148: Catherine The Great
$ws.Range("A148").Value = "Catherine The Great"
$ws.Range("B148").Value = "Robert Massie"
$ws.Range("C148").Value = (Get-Date -Year 2020 -Month 11 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D148").Value = (Get-Date -Year 2020 -Month 12 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E148").Value = "catherine the great;russia;18th century;autocracy;history;biography"
$ws.Range("F148").Value = "Audio"
$ws.Range("G148").Value = "23 Hours 52 Mins"
$ws.Range("H148").Value = 3
$ws.Range("I148").Value = $false

$ws.Range("A149").Select()
